$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45205 -> 45206) for every data row (rows 2 through 398).
$ws.Range("C2:C398").Value = 45206
